# Auto-generated script to update cryptos Price (D) and Volume(1h) (E) columns
# per the commit "Updated cryptos list on Thu Feb 15 23:13:09 UTC 2024 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue $ws "D2" '51.776.16'
Set-TextValue $ws "E2" '  -0.01%  '
Set-TextValue $ws "D3" '2.820.12'
Set-TextValue $ws "E3" '  +1.80%  '
Set-TextValue $ws "E4" '  -0.03%  '
Set-TextValue $ws "D5" '352.00'
Set-TextValue $ws "E5" '  +5.67%  '
Set-TextValue $ws "D6" '112.25'
Set-TextValue $ws "E6" '  -4.05%  '
Set-TextValue $ws "D7" '0.564'
Set-TextValue $ws "E7" '  +4.64%  '
Set-TextValue $ws "E8" '  -0.01%  '
Set-TextValue $ws "D9" '0.601'
Set-TextValue $ws "E9" '  +4.25%  '
Set-TextValue $ws "D10" '41.43'
Set-TextValue $ws "E10" '  -1.85%  '
Set-TextValue $ws "D11" '0.0851'
Set-TextValue $ws "E11" '  -0.92%  '
Set-TextValue $ws "E12" '  +1.13%  '
Set-TextValue $ws "D13" '19.88'
Set-TextValue $ws "E13" '  -2.04%  '
Set-TextValue $ws "D14" '7.75'
Set-TextValue $ws "E14" '  +1.28%  '
Set-TextValue $ws "D15" '3.255.06'
Set-TextValue $ws "E15" '  +1.54%  '
Set-TextValue $ws "D16" '2.823.97'
Set-TextValue $ws "E16" '  +1.22%  '
Set-TextValue $ws "D17" '0.885'
Set-TextValue $ws "E17" '  -0.68%  '
Set-TextValue $ws "D18" '51.597.74'
Set-TextValue $ws "E18" '  -0.18%  '
Set-TextValue $ws "D19" '7.41'
Set-TextValue $ws "E19" '  +8.03%  '
Set-TextValue $ws "D20" '3.19'
Set-TextValue $ws "E20" '  -2.16%  '
Set-TextValue $ws "D21" '13.40'
Set-TextValue $ws "E21" '  -1.38%  '
Set-TextValue $ws "D22" '0.0₃0992'
Set-TextValue $ws "E22" '  +1.66%  '
Set-TextValue $ws "D23" '269.50'
Set-TextValue $ws "E23" '  -3.65%  '
Set-TextValue $ws "D24" '69.66'
Set-TextValue $ws "E24" '  -0.53%  '
Set-TextValue $ws "D25" '2.76'
Set-TextValue $ws "E25" '  +2.24%  '
Set-TextValue $ws "D26" '26.76'
Set-TextValue $ws "E26" '  -0.32%  '
Set-TextValue $ws "E27" '  +0.10%  '
Set-TextValue $ws "D28" '10.31'
Set-TextValue $ws "E28" '  +0.93%  '
Set-TextValue $ws "D29" '2.25'
Set-TextValue $ws "E29" '  +0.83%  '
Set-TextValue $ws "E30" '  -1.93%  '
Set-TextValue $ws "D31" '50.66'
Set-TextValue $ws "E31" '  +1.08%  '
Set-TextValue $ws "D32" '33.89'
Set-TextValue $ws "E32" '  -3.52%  '
Set-TextValue $ws "D33" '0.0452'
Set-TextValue $ws "E33" '  +26.46%  '
Set-TextValue $ws "D34" '5.81'
Set-TextValue $ws "E34" '  +3.83%  '
Set-TextValue $ws "D35" '5.27'
Set-TextValue $ws "E35" '  +4.79%  '
Set-TextValue $ws "D36" '0.0821'
Set-TextValue $ws "E36" '  -0.09%  '
Set-TextValue $ws "E37" '  -0.13%  '
Set-TextValue $ws "D38" '2.06'
Set-TextValue $ws "E38" '  -1.85%  '
Set-TextValue $ws "D39" '3.22'
Set-TextValue $ws "E39" '  -0.98%  '
Set-TextValue $ws "D40" '18.14'
Set-TextValue $ws "E40" '  -5.60%  '
Set-TextValue $ws "D41" '23.80'
Set-TextValue $ws "E41" '  +1.96%  '
Set-TextValue $ws "E42" '  +1.82%  '
Set-TextValue $ws "D43" '125.97'
Set-TextValue $ws "E43" '  -0.95%  '
Set-TextValue $ws "D44" '2.52'
Set-TextValue $ws "E44" '  +2.53%  '
Set-TextValue $ws "E45" '  -1.26%  '
Set-TextValue $ws "D46" '2.080.19'
Set-TextValue $ws "E46" '  -0.37%  '
Set-TextValue $ws "D47" '3.32'
Set-TextValue $ws "E47" '  +0.03%  '
Set-TextValue $ws "E48" '  +1.65%  '
Set-TextValue $ws "D49" '5.65'
Set-TextValue $ws "E49" '  +2.03%  '
Set-TextValue $ws "E50" '  +5.98%  '
Set-TextValue $ws "D51" '60.44'
Set-TextValue $ws "E51" '  -0.41%  '
